$wb = $excel.ActiveWorkbook

# --- GAgg3: duration-bucket headers -> Korean short/medium/long labels ---
$wsGAgg3 = $wb.Worksheets.Item("GAgg3")
$wsGAgg3.Range("B1").Value = "단기"
$wsGAgg3.Range("C1").Value = "중기"
$wsGAgg3.Range("D1").Value = "장기"

# --- GAgg4: duration-bucket headers -> Korean labels + TIPS casing ---
$wsGAgg4 = $wb.Worksheets.Item("GAgg4")
$wsGAgg4.Range("B1").Value = "초단기"
$wsGAgg4.Range("C1").Value = "중단기"
$wsGAgg4.Range("D1").Value = "중기"
$wsGAgg4.Range("E1").Value = "장기"
$wsGAgg4.Range("F1").Value = "초장기"
$wsGAgg4.Range("G1").Value = "지방정부"
$wsGAgg4.Range("H1").Value = "TIPS"

# --- USIGSector: sector headers -> abbreviated labels ---
$wsUSIG = $wb.Worksheets.Item("USIGSector")
$wsUSIG.Range("B1").Value = "Cycl"
$wsUSIG.Range("C1").Value = "Def"
$wsUSIG.Range("D1").Value = "Ener"
$wsUSIG.Range("E1").Value = "Infra"
$wsUSIG.Range("F1").Value = "Util"
$wsUSIG.Range("G1").Value = "Bank"
$wsUSIG.Range("H1").Value = "Fin"

# --- restore/update the selection (active cell) on each sheet ---
$wsGAgg3.Activate()
$wsGAgg3.Range("E11").Select()

$wsGAgg4.Activate()
$wsGAgg4.Range("I15").Select()

$wsUSIG.Activate()
$wsUSIG.Range("G2").Select()

$wsGAgg1 = $wb.Worksheets.Item("GAgg1")
$wsGAgg1.Activate()
$wsGAgg1.Range("H15").Select()
